$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

$ws.Range("B2").Value = 94.893533687534145
$ws.Range("C2").Value = 93.419141172176538
$ws.Range("D2").Value = 93.318448052619615
$ws.Range("E2").Value = 94.292850339544714

$ws.Range("B3").Value = 93.690712321758312
$ws.Range("C3").Value = 93.786045175059073
$ws.Range("D3").Value = 91.747315679527389
$ws.Range("E3").Value = 95.751039117096752

$ws.Range("B1:E3").Select()
